$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to remain text so values like "1.00" or "0.390" are not
# reinterpreted as numbers and lose their exact formatting, matching the source inlineStr cells.
$priceCells = @("D2", "D3", "D5", "D8", "D9", "D12", "D13", "D14", "D16", "D17", "D18", "D20", "D22", "D23", "D26", "D27", "D29", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D48", "D49", "D51")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "37.543.39"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "2.082.49"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "233.12"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "57.89"
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").Value = "0.390"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("D12").Value = "15.06"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").Value = "2.389.64"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "21.10"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "5.36"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "2.079.78"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "37.513.09"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("D20").Value = "70.76"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "228.89"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").Value = "9.69"
$ws.Range("E26").Value = "  +7.30%  "
$ws.Range("D27").Value = "170.45"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -3.99%  "
$ws.Range("D29").Value = "19.70"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "5.31"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "0.0232"
$ws.Range("E40").Value = "  +7.65%  "
$ws.Range("D41").Value = "100.35"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "0.0957"
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "2.91"
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("E44").Value = "  +3.65%  "
$ws.Range("D45").Value = "16.86"
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("D46").Value = "1.462.39"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "3.96"
$ws.Range("E48").Value = "  -5.89%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "7.26"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").Value = "2.273.57"
$ws.Range("E51").Value = "  +0.18%  "
